$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- 1. Notes page: replace the 4x "Notes" paragraphs + trailing empty
#        paragraph with an empty paragraph, a "Today is 09/21" paragraph,
#        and a trailing empty paragraph.
$notesPage = $s.NotesPage
$notesShape = $notesPage.Shapes.Item(2)
$notesShape.TextFrame.TextRange.Text = "`nToday is 09/21`n"

# --- 2. Table formatting: drop the special "first row" banding, add
#        centered alignment to the top-left header cell, shrink every
#        run to 8pt and collapse every row (and the frame) to minimum
#        height.
$tblShape = $s.Shapes.Item(3)
$tbl = $tblShape.Table

$tbl.FirstRow = $false

$rowCount = $tbl.Rows.Count
$colCount = $tbl.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Shape.TextFrame.TextRange.Font.Size = 8
    }
    $tbl.Rows.Item($r).Height = 0
}

$tbl.Cell(1, 1).Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$tblShape.Height = 0
